# Generate Report for Handback
# Update the handoff/handback timestamp values that were regenerated when
# the handback report was (re)generated for the e72579d7-... file.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# "Latest HO Xliff Generate Date" for e72579d7-...md
# This value is shared between Overview!G4 and de-de!H4 (same underlying
# shared string), so both must be updated together.
$wsOverview.Range("G4").Value = "2016-08-17 00:42:47"
$wsDeDe.Range("H4").Value = "2016-08-17 00:42:47"

# zh-cn row for e72579d7-...: Correspond Handoff Datetime / Correspond Handback DateTime
$wsZhCn.Range("H4").Value = "2016-08-17 00:42:41"
$wsZhCn.Range("K4").Value = "2016-08-17 00:42:59"

# de-de row for e72579d7-...: Correspond Handback DateTime
$wsDeDe.Range("K4").Value = "2016-08-17 00:43:11"
